# Master v2 - WIP
# Update abbreviation definitions on the "llm" sheet and drop the trailing
# two rows (UPPA / SCM) that were appended in error.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("llm")

$ws.Range("B2").Value = "Ecole Polytechnique de l'Ouest de la Guinée"
$ws.Range("B6").Value = "Institut des Stratégies et Techniques de Communication"
$ws.Range("B7").Value = "Concours d'Accès au Corps des Formateurs de Personnels"
$ws.Range("B14").Value = "Institut National de Perfectionnement et d'Éducateurs"
$ws.Range("B17").Value = "École Supérieure des Sciences et Langues pour l'Ingénieur de Lille"
$ws.Range("B19").Value = "..."
$ws.Range("B39").Value = "Licence Langues, Littératures et Civilisations Étrangères et Régionales"
$ws.Range("B42").Value = "École de Management de la Chambre de Commerce et d'Industrie de Paris"
$ws.Range("B57").Value = "Université de Science et Technologie de Hanoï"
$ws.Range("B67").Value = "..."
$ws.Range("B76").Value = "Sciences Humaines et Sociales, Philosophie, Sociologie"
$ws.Range("B178").Value = "Diplôme Universitaire de Musicien Intervenant"
$ws.Range("B184").Value = "Sciences Humaines et Sociales, Philosophie"
$ws.Range("B196").Value = "Institut Supérieur de Formation de l'Enseignement Catholique"

$ws.Range("B205").Value = ""

$ws.Rows.Item(207).Delete()
$ws.Rows.Item(206).Delete()
